$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-history snapshot column was added just before the existing
# "nom" / "url_produit" columns. Inserting a whole column at BL pushes
# "nom" from BL to BM and "url_produit" from BM to BN for every row,
# exactly mirroring the diff's column shift.
$ws.Columns("BL:BL").Insert()

# Header for the freshly inserted column: the newest scrape timestamp.
$ws.Range("BL1").Value = "2026-01-30 14:26:31"

# For the real product rows (2-80), the new timestamp column repeats the
# most recently recorded price, i.e. the same value already sitting in
# column BK for that row.
for ($row = 2; $row -le 80; $row++) {
    $ws.Cells.Item($row, 63).Copy($ws.Cells.Item($row, 64))
}
